# Updates the cryptos list: refresh the Price (D) / Volume(1h) (E) columns
# for each data row, and swap the Filecoin / ImmutableX rows (32 & 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $d, $e) {
    if ($d -ne $null) {
        $ws.Cells.Item($row, 4).Value = $d
    }
    if ($e -ne $null) {
        $ws.Cells.Item($row, 5).Value = "  $e  "
    }
}

Set-Row 2  "20.584.05"   "+0.40%"
Set-Row 3  "1.479.75"    "+0.62%"
Set-Row 4  $null          "+0.26%"
Set-Row 5  "0.9709"      "+1.37%"
Set-Row 6  "279.26"      "-0.72%"
Set-Row 7  "0.3661"      "-1.12%"
Set-Row 8  "0.3079"      "-3.10%"
Set-Row 9  "40.00"       "-4.28%"
Set-Row 10 "1.064"       "+0.64%"
Set-Row 11 "0.06677"     "+0.03%"
Set-Row 12 $null          "+0.02%"
Set-Row 13 "5.527"       "-1.35%"
Set-Row 14 "18.05"       "-0.95%"
Set-Row 15 "6.217"       $null
Set-Row 16 "0.9721"      "+1.46%"
Set-Row 17 "0.00001029"  "-0.64%"
Set-Row 18 "1.478.64"    "+0.18%"
Set-Row 19 "0.05939"     "+3.79%"
Set-Row 20 "69.69"       "-3.45%"
Set-Row 21 "5.499"       "-2.99%"
Set-Row 22 "14.52"       "-1.11%"
Set-Row 23 $null          "-1.37%"
Set-Row 24 "2.261"       "+0.32%"
Set-Row 25 "20.642.92"   "-0.24%"
Set-Row 26 "141.95"      "+2.86%"
Set-Row 27 "2.130"       "-7.01%"
Set-Row 28 "17.30"       "-1.52%"
Set-Row 29 "1.639.09"    "+0.06%"
Set-Row 30 "114.14"      "+0.34%"
Set-Row 31 "3.923"       "-0.89%"

# Rows 32 & 33 swap: Filecoin <-> ImmutableX (coin name, link, price, volume)
$ws.Cells.Item(32, 2).Value = "ImmutableX"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(32, 4).Value = "0.8239"
$ws.Cells.Item(32, 5).Value = "  -0.93%  "

$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(33, 4).Value = "5.017"
$ws.Cells.Item(33, 5).Value = "  -5.62%  "

Set-Row 34 "0.07988"     "+1.98%"
Set-Row 35 "1.534"       "-5.73%"
Set-Row 36 "1.208"       "+7.67%"
Set-Row 37 "0.05790"     "-4.06%"
Set-Row 38 "4.732"       "-3.71%"
Set-Row 39 "0.9716"      "+0.12%"
Set-Row 40 $null          "-1.22%"
Set-Row 41 "10.48"       "-1.47%"
Set-Row 42 "7.658"       "+4.56%"
Set-Row 43 "0.1879"      "-0.53%"
Set-Row 44 "0.5304"      "-1.94%"
Set-Row 45 "3.533"       "-1.60%"
Set-Row 46 "12.27"       "-1.16%"
Set-Row 47 "118.50"      "-2.85%"
Set-Row 48 "0.5198"      "-2.11%"
Set-Row 49 "1.808"       "-0.55%"
Set-Row 50 $null          "+0.83%"
Set-Row 51 "0.9946"      "+0.14%"
